$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of writes matters: new shared-string entries get appended in the
# order they're first introduced, so these are sequenced to match the
# target workbook's shared-strings table ordering.

# --- D6: new text, existing style (s=26) kept as-is ---
$ws.Range("D6").Value  = "Intro to R part 4"

# --- C7: new text, existing style (s=31) kept as-is ---
$ws.Range("C7").Value  = "Stats Lab 1: Camera data organization and exploration"

# --- E7: new text, existing style (s=25) kept as-is ---
$ws.Range("E7").Value  = "Intro to R part 5"

# --- D9: new text, existing style (s=25) kept as-is ---
$ws.Range("D9").Value  = "Data Visualization Part 1"

# --- E9 needs both new text and a style change (s=18 -> s=25, matching D9's style) ---
$ws.Range("D9").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "Data Visualization Part 2"

# --- E10 needs both new text and a style change (s=18 -> s=25, matching D9's style) ---
$ws.Range("D9").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "Stats Lab 2: Generalized linear modelling I"

# --- D12: new text, existing style (s=25) kept as-is ---
$ws.Range("D12").Value = "Stats Lab 3: Modelling Your Data I"

# --- D13: new text, existing style (s=25) kept as-is ---
$ws.Range("D13").Value = "Stats Lab 4: Modelling Your Data II"

# --- D15: new text, existing style (s=25) kept as-is ---
$ws.Range("D15").Value = "Stats Lab 5: Model validation and presentation"

# --- F12 is a brand-new cell; copy F10's format (s=28) then set its value ---
$ws.Range("F10").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = "Stats Assignment 6 DUE"

# --- E6: new text, existing style (s=25) kept as-is ---
$ws.Range("E6").Value  = "Intro to R part 4 cont…"

# --- Selection moved from D10 to E6 ---
$ws.Range("E6").Select()
